$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Sunday, Jan 15", "11:45 AM", "FR2137", "London", "(STN)", "Lauda Europe ", "A320", "(9H-LOT)", "11:52 AM", "0 hours, 7 minutes"),
    @("Sunday, Jan 15", "12:05 PM", "LO3810", "Warsaw", "(WAW)", "LOT ", "E195", "(SP-LNN)", "12:08 PM", "0 hours, 3 minutes"),
    @("Sunday, Jan 15", "1:30 PM", "5Y9445", "Hong Kong", "(HKG)", "Atlas Air ", "B744", "(N485MC)", "1:50 PM", "0 hours, 20 minutes"),
    @("Sunday, Jan 15", "2:55 PM", "LO3802", "Warsaw", "(WAW)", "LOT ", "E190", "(SP-LMA)", "3:07 PM", "0 hours, 12 minutes")
)

$startRow = 78
$startNumber = 77

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = [double]($startNumber + $i)
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $ws.Cells.Item($r, 9).Value = $row[7]
    $ws.Cells.Item($r, 10).Value = $row[8]
    $ws.Cells.Item($r, 12).Value = $row[9]
}
